$d = $word.ActiveDocument

# "...total hits from each team, total outs, as well..." -> "...total runs, as well..."
$d.Content.Find.Execute("total outs", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "total runs", 2) | Out-Null

# "...the type of throw the pitcher used, the speed..." -> "...the number of strikeouts from the pitcher, the speed..."
$d.Content.Find.Execute("the type of throw the pitcher used, the speed", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "the number of strikeouts from the pitcher, the speed", 2) | Out-Null

# Note: "interesting" already sits at the end of the run immediately before the page break, and the
# following run already begins with "information that will be presented..." -- the concatenated text
# "...compute interesting information..." is unchanged, so no edit is needed there. Leaving it alone
# also avoids touching the run that carries <w:lastRenderedPageBreak/>.
